$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "protected and private inheritance"
$ws.Range("C20").Value = "name=prasuna"

# Row 21 - first hyperlink (creates shared string "mail=@gmail.com")
$ws.Range("C21").Value = "mail=@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:prasuna@gmail.com")

# Finish row 20
$ws.Range("D20").Value = "protected inheritace"
$ws.Range("E20").Value = "protected inheritace"

# Finish row 21
$ws.Range("D21").Value = "name=prasuna"
$ws.Range("E21").Value = "name=prasuna"

# Row 23 text cells (creates "private inheritance")
$ws.Range("D23").Value = "private inheritance"
$ws.Range("E23").Value = "private inheritance"

# Row 22 text cell (creates "name=a125")
$ws.Range("C22").Value = "name=a125"

# Row 24
$ws.Range("D24").Value = "name=a125"
$ws.Range("E24").Value = "name=a125"

# Row 22 result column
$ws.Range("F22").Value = "PASS"

# Remaining hyperlinks (reuse "mail=@gmail.com")
$ws.Range("C23").Value = "mail=@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C23"), "mailto:a125@gmail.com")

$ws.Range("D22").Value = "mail=@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D22"), "mailto:a125@gmail.com")

$ws.Range("D25").Value = "mail=@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D25"), "mailto:a125@gmail.com")

$ws.Range("E22").Value = "mail=@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E22"), "mailto:a125@gmail.com")

$ws.Range("E25").Value = "mail=@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E25"), "mailto:a125@gmail.com")

# Update the view selection to match the new active cell state
$ws.Range("F23").Select() | Out-Null
